# Generate Report for Handoff
# Adds two new rows (for newly-tracked files 4217bb01-... and 53421fbb-...)
# to the Overview, zh-cn and de-de sheets of the localization-status workbook.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

$uuid1 = "4217bb01-cd18-4092-bbc2-c7c7c4cc7898"
$uuid2 = "53421fbb-37cb-4242-906f-0b14c04bc966"

$hash1 = "e32bd89353f79594a0303fafe9dbadd98a86a373"
$hash2 = "8d4f7b293469a8cc288318eb43d447c3f82b9aea"

$status = "Ready for handoff"

# ---------------------------------------------------------------------------
# Overview sheet: File Name | zh-cn | de-de | Latest Handoff Date
# ---------------------------------------------------------------------------

$ws1.Hyperlinks.Add($ws1.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/31714bb9e92751df0a81f02819205b3a53dc4480/e2e/$uuid1.md", "", "", "$uuid1.md") | Out-Null
$ws1.Range("B4").Value2 = $status
$ws1.Range("C4").Value2 = $status
$ws1.Range("D4").Value2 = "2016-29-20 10:29:36"

$ws1.Hyperlinks.Add($ws1.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/3a37d86fbaf657e7767975ee7e956732fb6983a9/e2e/$uuid2.md", "", "", "$uuid2.md") | Out-Null
$ws1.Range("B5").Value2 = $status
$ws1.Range("C5").Value2 = $status
$ws1.Range("D5").Value2 = "2016-29-20 10:29:36"

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------

$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/31714bb9e92751df0a81f02819205b3a53dc4480/e2e/$uuid1.md", "", "", "$uuid1.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("B4"), "https://github.com/OpenLocalizationTest/oltest/blob/31714bb9e92751df0a81f02819205b3a53dc4480/e2e/$uuid1.md", "", "", ".md") | Out-Null
$ws2.Range("C4").Value2 = $status
$ws2.Hyperlinks.Add($ws2.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/195db71da40f6117af98513b166ea973354bf675/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$uuid1.$hash1.zh-cn.xlf", "", "", "$uuid1.$hash1.zh-cn.xlf") | Out-Null
$ws2.Range("E4").Value2 = "2016-03-20 10:29:33"
$ws2.Range("H4").Value2 = "0001-01-01 00:00:00"
$ws2.Range("I4").Value2 = "Include"

$ws2.Hyperlinks.Add($ws2.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/3a37d86fbaf657e7767975ee7e956732fb6983a9/e2e/$uuid2.md", "", "", "$uuid2.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("B5"), "https://github.com/OpenLocalizationTest/oltest/blob/3a37d86fbaf657e7767975ee7e956732fb6983a9/e2e/$uuid2.md", "", "", ".md") | Out-Null
$ws2.Range("C5").Value2 = $status
$ws2.Hyperlinks.Add($ws2.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/bf2e2db6ff5bb75b2a0b11d1b7bbe7bf3ca0f4ca/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$uuid2.$hash2.zh-cn.xlf", "", "", "$uuid2.$hash2.zh-cn.xlf") | Out-Null
$ws2.Range("E5").Value2 = "2016-03-20 10:29:33"
$ws2.Range("H5").Value2 = "0001-01-01 00:00:00"
$ws2.Range("I5").Value2 = "Include"

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------

$ws3.Hyperlinks.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/31714bb9e92751df0a81f02819205b3a53dc4480/e2e/$uuid1.md", "", "", "$uuid1.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("B4"), "https://github.com/OpenLocalizationTest/oltest/blob/31714bb9e92751df0a81f02819205b3a53dc4480/e2e/$uuid1.md", "", "", ".md") | Out-Null
$ws3.Range("C4").Value2 = $status
$ws3.Hyperlinks.Add($ws3.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d6b1b5f10f95c6adb466904d5e89841c02832eed/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$uuid1.$hash1.de-de.xlf", "", "", "$uuid1.$hash1.de-de.xlf") | Out-Null
$ws3.Range("E4").Value2 = "2016-03-20 10:29:36"
$ws3.Range("H4").Value2 = "0001-01-01 00:00:00"
$ws3.Range("I4").Value2 = "Include"

$ws3.Hyperlinks.Add($ws3.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/3a37d86fbaf657e7767975ee7e956732fb6983a9/e2e/$uuid2.md", "", "", "$uuid2.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("B5"), "https://github.com/OpenLocalizationTest/oltest/blob/3a37d86fbaf657e7767975ee7e956732fb6983a9/e2e/$uuid2.md", "", "", ".md") | Out-Null
$ws3.Range("C5").Value2 = $status
$ws3.Hyperlinks.Add($ws3.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/15c4fef1e41cfd89cc7c90c2e2f8a2a57d7a4b0a/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$uuid2.$hash2.de-de.xlf", "", "", "$uuid2.$hash2.de-de.xlf") | Out-Null
$ws3.Range("E5").Value2 = "2016-03-20 10:29:36"
$ws3.Range("H5").Value2 = "0001-01-01 00:00:00"
$ws3.Range("I5").Value2 = "Include"

Write-Host "Done applying handback report updates"
